$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap I1/J1 ---
$ws.Range("I1").Value = 'SinglePrice'
$ws.Range("J1").Value = 'EnglishDescription'

# --- Barcode column (H) must stay text; set format before values ---
$ws.Range("H2:H28").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 30
$ws.Range("H2").Value = '6931286064278'
$ws.Range("I2").Value = 2.07
$ws.Range("J2").Value = 'Wolong- (Wuzi Scent) Old Stove Pot 138g'

# Row 3
$ws.Range("A3").Value = 30
$ws.Range("H3").Value = '6931286065268'
$ws.Range("I3").Value = 2.07
$ws.Range("J3").Value = 'Wolong- (Ecstasy Spicy) Old Stove Pot 128g'

# Row 4
$ws.Range("A4").Value = 10
$ws.Range("H4").Value = '6971657871364'
$ws.Range("I4").Value = 2.1
$ws.Range("J4").Value = 'San Yuanze-Hawthorn pulp Beijing cakes 160g'

# Row 5
$ws.Range("A5").Value = 15
$ws.Range("H5").Value = '6902083815103'
$ws.Range("I5").Value = 1.54
$ws.Range("J5").Value = 'Wahaha- (sugar-free) Dahongpao oolong tea 500ml'

# Row 6
$ws.Range("A6").Value = 15
$ws.Range("H6").Value = '6902083815110'
$ws.Range("I6").Value = 1.54
$ws.Range("J6").Value = 'Wahaha- (sugar-free) Zhengshan small breeding black tea 500ml'

# Row 7
$ws.Range("A7").Value = 15
$ws.Range("H7").Value = '6902083815134'
$ws.Range("I7").Value = 1.54
$ws.Range("J7").Value = 'Wahaha- (Sugar-free) Pu''er tea 500ml'

# Row 8
$ws.Range("A8").Value = 15
$ws.Range("H8").Value = '6902083815127'
$ws.Range("I8").Value = 1.54
$ws.Range("J8").Value = 'Wahaha- (sugar-free) jasmine tea 500ml'

# Row 9
$ws.Range("A9").Value = 10
$ws.Range("H9").Value = '6926896705277'
$ws.Range("I9").Value = 6.72
$ws.Range("J9").Value = 'Jixiangju-Crispy double bamboo shoot 106g'

# Row 10
$ws.Range("A10").Value = 10
$ws.Range("H10").Value = '6909003200117'
$ws.Range("I10").Value = 1.13
$ws.Range("J10").Value = 'Cuckoo City-150g of tempeh'

# Row 11
$ws.Range("A11").Value = 10
$ws.Range("H11").Value = '6907592001429'
$ws.Range("I11").Value = 4.25
$ws.Range("J11").Value = 'Wang Zhihe-dried yellow sauce 180g'

# Row 12
$ws.Range("A12").Value = 24
$ws.Range("H12").Value = '6928002376630'
$ws.Range("I12").Value = 2.41
$ws.Range("J12").Value = 'Gan Juyuan-pure brown sugar 300g'

# Row 13
$ws.Range("A13").Value = 15
$ws.Range("H13").Value = '6922507808597'
$ws.Range("I13").Value = 2.84
$ws.Range("J13").Value = 'Chen Keming- (Egg) Refined Noodle 1000g'

# Row 14
$ws.Range("A14").Value = 15
$ws.Range("H14").Value = '6922507808580'
$ws.Range("I14").Value = 2.84
$ws.Range("J14").Value = 'Chen Keming- (Jin Dao) Refined Noodle 1000g'

# Row 15
$ws.Range("A15").Value = 50
$ws.Range("H15").Value = '6920404395981'
$ws.Range("I15").Value = 0.27
$ws.Range("J15").Value = 'Nongfu Villa-California Ximei 70g'

# Row 16
$ws.Range("A16").Value = 24
$ws.Range("H16").Value = '6932419100474'
$ws.Range("I16").Value = 8.75
$ws.Range("J16").Value = 'Shuangrong- (Rose Black Sugar) Guanyin Pavilion 220g'

# Row 17
$ws.Range("A17").Value = 30
$ws.Range("C17").Value = '伊藤园-原味烏龍茶 500ml'
$ws.Range("H17").Value = '6958959708208'
$ws.Range("I17").Value = 0.91
$ws.Range("J17").Value = 'Ito Garden-Original oolong tea 500ml'

# Row 18
$ws.Range("A18").Value = 30
$ws.Range("C18").Value = '伊藤园-無糖冷萃綠茶 600ml'
$ws.Range("H18").Value = '6958959708246'
$ws.Range("I18").Value = 0.92
$ws.Range("J18").Value = 'Ito Garden-sugar-free cold green tea 600ml'

# Row 19
$ws.Range("A19").Value = 30
$ws.Range("C19").Value = '白家-酸湯肥牛調料 160g'
$ws.Range("H19").Value = '6974142070525'
$ws.Range("I19").Value = 2.42
$ws.Range("J19").Value = 'Baijia-Sour soup Fat beef seasoning 160g'

# Row 20
$ws.Range("A20").Value = 30
$ws.Range("C20").Value = '白家-老鹵汁調料 120g'
$ws.Range("H20").Value = '6926410335416'
$ws.Range("I20").Value = 1.89
$ws.Range("J20").Value = 'Baijia-Old marinade seasoning 120g'

# Row 21
$ws.Range("A21").Value = 15
$ws.Range("C21").Value = '道地-蘋果綠茶 500ml'
$ws.Range("H21").Value = '6976901700001'
$ws.Range("I21").Value = 1.85
$ws.Range("J21").Value = 'Authentic-apple green tea 500ml'

# Row 22
$ws.Range("A22").Value = 15
$ws.Range("C22").Value = '道地-巨峰紅茶 500ml'
$ws.Range("H22").Value = '6976901700018'
$ws.Range("I22").Value = 1.85
$ws.Range("J22").Value = 'Authentic-Jufeng Black Tea 500ml'

# Row 23
$ws.Range("A23").Value = 12
$ws.Range("C23").Value = 'Calpis 原味可爾必思 990ml'
$ws.Range("H23").Value = '4714947001131'
$ws.Range("I23").Value = 2.71
$ws.Range("J23").Value = 'Calpis Original Calle Betis 990ml'

# Row 24
$ws.Range("A24").Value = 24
$ws.Range("C24").Value = 'Calpis 原味可爾必思 500ml'
$ws.Range("H24").Value = '4714947000134'
$ws.Range("I24").Value = 1.59
$ws.Range("J24").Value = 'Calpis original flavor Kerbis Betis 500ml'

# Row 25
$ws.Range("A25").Value = 24
$ws.Range("C25").Value = 'Calpis 葡萄可爾必思 500ml'
$ws.Range("H25").Value = '4714947013615'
$ws.Range("I25").Value = 1.59
$ws.Range("J25").Value = 'Calpis Grape Calbin 500ml'

# Row 26
$ws.Range("A26").Value = 24
$ws.Range("C26").Value = '成央記-海南老椰汁 420ml'
$ws.Range("H26").Value = '6972698570087'
$ws.Range("I26").Value = 3.85
$ws.Range("J26").Value = 'Cheng Yang Ji-Hainan Old Coconut Sprathy 420ml'

# Row 27
$ws.Range("A27").Value = 30
$ws.Range("C27").Value = '道地-烏龍茶 500ml'
$ws.Range("H27").Value = '6930599700026'
$ws.Range("I27").Value = 0.92
$ws.Range("J27").Value = 'Authentic-oolong tea 500ml'

# Row 28
$ws.Range("A28").Value = 30
$ws.Range("C28").Value = '道地-解。綠茶 500ml'
$ws.Range("H28").Value = '6930599700101'
$ws.Range("I28").Value = 0.92
$ws.Range("J28").Value = 'Authentic-solution.Green Tea 500ml'
